$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its original text formatting so that
# numeric-looking strings (e.g. "233.12") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.045.32"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "1.816.61"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "233.12"
$ws.Range("E5").Value = "  +3.00%  "

$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").Value = "40.39"
$ws.Range("E8").Value = "  -7.86%  "

$ws.Range("E9").Value = "  +10.35%  "

$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Value = "2.078.45"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").Value = "1.820.05"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").Value = "11.08"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").Value = "0.662"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("D17").Value = "34.999.48"
$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("D18").Value = "69.53"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").Value = "238.96"
$ws.Range("E20").Value = "  -1.39%  "

$ws.Range("D21").Value = "11.84"
$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("D22").Value = "4.65"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("E24").Value = "  +3.19%  "

$ws.Range("D25").Value = "172.44"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("E28").Value = "  -1.17%  "

$ws.Range("E29").Value = "  +29.53%  "

$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").Value = "3.339.69"
$ws.Range("E31").Value = "  +37.45%  "

$ws.Range("D32").Value = "0.0556"
$ws.Range("E32").Value = "  +6.17%  "

$ws.Range("D33").Value = "3.94"
$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("E35").Value = "  -3.60%  "

$ws.Range("D36").Value = "93.09"
$ws.Range("E36").Value = "  +3.74%  "

$ws.Range("E37").Value = "  +6.93%  "

$ws.Range("E38").Value = "  +2.77%  "

$ws.Range("D39").Value = "0.0193"
$ws.Range("E39").Value = "  +0.94%  "

$ws.Range("D40").Value = "1.28"
$ws.Range("E40").Value = "  +4.05%  "

$ws.Range("D41").Value = "1.306.08"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").Value = "0.986"
$ws.Range("E42").Value = "  +1.88%  "

$ws.Range("E43").Value = "  -3.55%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "14.59"
$ws.Range("E44").Value = "  -4.39%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  +1.33%  "

$ws.Range("E46").Value = "  -1.68%  "

$ws.Range("D47").Value = "6.34"
$ws.Range("E47").Value = "  +6.54%  "

$ws.Range("D48").Value = "0.0510"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("D49").Value = "1.993.61"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").Value = "0.0645"
$ws.Range("E51").Value = "  +5.09%  "
